# This workbook was re-saved by a newer version of Excel (the author simply
# opened the file and saved it again after finishing the project). The
# functional/visible changes captured in the diff are:
#
#   1. Column A ("Question") was widened to a custom width of 81 characters.
#   2. Column B ("Reponse Attendue") was widened to a custom width of
#      ~228.57 characters (an autofit-style width to accommodate the long
#      answer text).
#   3. The active selection was left on cell A18 when the file was saved.
#   4. `fullCalcOnLoad` was dropped from the workbook (no longer forces a
#      full recalculation on open) -- calculation settings otherwise
#      unchanged.
#
# (Everything else in the raw XML diff -- extra mc/x14ac/xr namespace
# declarations, the xr:uid/revisionPtr GUIDs, the hard-coded absPath to the
# author's desktop, the window position/size, and the theme's cosmetic
# "2007 - 2010" renaming -- are metadata Excel stamps automatically onto a
# file purely because a newer build opened and re-saved it; none of it is
# reachable through a scripted object-model edit, so it is intentionally
# left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (Question) and column B (Reponse Attendue) to fit their
# (long) text content.
$ws.Columns.Item(1).ColumnWidth = 80.2
$ws.Columns.Item(2).ColumnWidth = 227.6

# Leave the selection on A18, matching the saved cursor position.
$ws.Range("A18").Select() | Out-Null
